$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price / volume data refreshed by the scraper run. Numeric-looking values
# are written with a leading apostrophe so Excel keeps them as literal text (the
# sheet stores prices like "26.170.16" / "19.15" and padded percents as strings,
# not as numbers), then the style is reset to Normal so no stray number-format
# gets attached to the cell (keeping cells unstyled, like the original file).
$ws.Range("D2").Value = "'26.170.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.70%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.585.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'211.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.94%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.29%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.59%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.03%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.21%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.809.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.40%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.587.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.50%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.31%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'63.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.14%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.185.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.65%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.60%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'214.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.20%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -3.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.03%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.73%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.33%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.67%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'144.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.68%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'6.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.26%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.26%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.37%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.91%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.59%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D33").Value = "'1.406.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +7.89%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.60%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.51%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.36%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.35%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.85%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +3.99%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.943"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -14.64%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.766"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.36%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.720.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.47%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'60.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.85%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'85.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.84%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.14%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Cronos"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0500"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Algorand"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0970"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.40%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.25%  "
$ws.Range("E51").Style = "Normal"
